$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap the F:V (match detail) content between paired fixtures that were
# originally recorded out of kickoff-time order. A:E (index/country/
# tournament/season/date) are identical within each pair, so only F:V move.
$tmp10 = $ws.Range("F10:V10").Value2
$tmp11 = $ws.Range("F11:V11").Value2
$ws.Range("F10:V10").Value2 = $tmp11
$ws.Range("F11:V11").Value2 = $tmp10

$tmp40 = $ws.Range("F40:V40").Value2
$tmp41 = $ws.Range("F41:V41").Value2
$ws.Range("F40:V40").Value2 = $tmp41
$ws.Range("F41:V41").Value2 = $tmp40

$tmp48 = $ws.Range("F48:V48").Value2
$tmp49 = $ws.Range("F49:V49").Value2
$ws.Range("F48:V48").Value2 = $tmp49
$ws.Range("F49:V49").Value2 = $tmp48

$tmp64 = $ws.Range("F64:V64").Value2
$tmp65 = $ws.Range("F65:V65").Value2
$ws.Range("F64:V64").Value2 = $tmp65
$ws.Range("F65:V65").Value2 = $tmp64

$tmp75 = $ws.Range("F75:V75").Value2
$tmp76 = $ws.Range("F76:V76").Value2
$ws.Range("F75:V75").Value2 = $tmp76
$ws.Range("F76:V76").Value2 = $tmp75

# --- Append 5 new fixtures (rows 80-84) for matchday played 07/01/2024.
$ws.Range("A79:V79").Copy()
$ws.Range("A80:V84").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$row80 = New-Object 'object[,]' 1,22
$row80[0,0] = 79
$row80[0,1] = "malta"
$row80[0,2] = "premier-league"
$row80[0,3] = "2023-2024"
$row80[0,4] = 45298.45833333334
$row80[0,5] = "Sliema"
$row80[0,6] = 2
$row80[0,7] = "Gudja"
$row80[0,8] = 0
$row80[0,9] = 1.44
$row80[0,10] = "05/01/2024 23:12"
$row80[0,11] = 1.33
$row80[0,12] = "07/01/2024 10:59"
$row80[0,13] = 3.84
$row80[0,14] = "05/01/2024 23:12"
$row80[0,15] = 4.55
$row80[0,16] = "07/01/2024 10:59"
$row80[0,17] = 6.66
$row80[0,18] = "05/01/2024 23:12"
$row80[0,19] = 10.21
$row80[0,20] = "07/01/2024 10:59"
$row80[0,21] = "https://www.betexplorer.com/football/malta/premier-league/sliema-gudja/feRoXhSJ/"
$ws.Range("A80:V80").Value2 = $row80

$row81 = New-Object 'object[,]' 1,22
$row81[0,0] = 80
$row81[0,1] = "malta"
$row81[0,2] = "premier-league"
$row81[0,3] = "2023-2024"
$row81[0,4] = 45298.58333333334
$row81[0,5] = "Marsaxlokk"
$row81[0,6] = 1
$row81[0,7] = "Hibernians"
$row81[0,8] = 1
$row81[0,9] = 2.64
$row81[0,10] = "06/01/2024 02:12"
$row81[0,11] = 2.65
$row81[0,12] = "07/01/2024 13:55"
$row81[0,13] = 2.89
$row81[0,14] = "06/01/2024 02:12"
$row81[0,15] = 3.43
$row81[0,16] = "07/01/2024 13:56"
$row81[0,17] = 2.63
$row81[0,18] = "06/01/2024 02:12"
$row81[0,19] = 2.49
$row81[0,20] = "07/01/2024 13:52"
$row81[0,21] = "https://www.betexplorer.com/football/malta/premier-league/marsaxlokk-hibernians/vT2IRWks/"
$ws.Range("A81:V81").Value2 = $row81

$row82 = New-Object 'object[,]' 1,22
$row82[0,0] = 81
$row82[0,1] = "malta"
$row82[0,2] = "premier-league"
$row82[0,3] = "2023-2024"
$row82[0,4] = 45298.58333333334
$row82[0,5] = "Sirens"
$row82[0,6] = 1
$row82[0,7] = "Birkirkara"
$row82[0,8] = 2
$row82[0,9] = 5.73
$row82[0,10] = "06/01/2024 02:12"
$row82[0,11] = 6.5
$row82[0,12] = "07/01/2024 12:05"
$row82[0,13] = 3.78
$row82[0,14] = "06/01/2024 02:12"
$row82[0,15] = 4.08
$row82[0,16] = "07/01/2024 12:05"
$row82[0,17] = 1.5
$row82[0,18] = "06/01/2024 02:12"
$row82[0,19] = 1.49
$row82[0,20] = "07/01/2024 12:05"
$row82[0,21] = "https://www.betexplorer.com/football/malta/premier-league/sirens-birkirkara/YoQkWCsQ/"
$ws.Range("A82:V82").Value2 = $row82

$row83 = New-Object 'object[,]' 1,22
$row83[0,0] = 82
$row83[0,1] = "malta"
$row83[0,2] = "premier-league"
$row83[0,3] = "2023-2024"
$row83[0,4] = 45298.67708333334
$row83[0,5] = "Mosta"
$row83[0,6] = 0
$row83[0,7] = "Balzan"
$row83[0,8] = 2
$row83[0,9] = 2.56
$row83[0,10] = "06/01/2024 04:42"
$row83[0,11] = 2.91
$row83[0,12] = "07/01/2024 16:12"
$row83[0,13] = 3.19
$row83[0,14] = "06/01/2024 04:42"
$row83[0,15] = 3.08
$row83[0,16] = "07/01/2024 16:06"
$row83[0,17] = 2.5
$row83[0,18] = "06/01/2024 04:42"
$row83[0,19] = 2.48
$row83[0,20] = "07/01/2024 16:12"
$row83[0,21] = "https://www.betexplorer.com/football/malta/premier-league/mosta-fc-balzan-fc/ng6QPAKg/"
$ws.Range("A83:V83").Value2 = $row83

$row84 = New-Object 'object[,]' 1,22
$row84[0,0] = 83
$row84[0,1] = "malta"
$row84[0,2] = "premier-league"
$row84[0,3] = "2023-2024"
$row84[0,4] = 45298.67708333334
$row84[0,5] = "Valletta"
$row84[0,6] = 1
$row84[0,7] = "Naxxar"
$row84[0,8] = 2
$row84[0,9] = 1.69
$row84[0,10] = "06/01/2024 04:42"
$row84[0,11] = 1.79
$row84[0,12] = "07/01/2024 16:06"
$row84[0,13] = 3.49
$row84[0,14] = "06/01/2024 04:42"
$row84[0,15] = 3.5
$row84[0,16] = "07/01/2024 16:06"
$row84[0,17] = 4.37
$row84[0,18] = "06/01/2024 04:42"
$row84[0,19] = 4.38
$row84[0,20] = "07/01/2024 16:06"
$row84[0,21] = "https://www.betexplorer.com/football/malta/premier-league/valletta-naxxar-lions/hW6MQj5m/"
$ws.Range("A84:V84").Value2 = $row84

